# Femacal de La Calera - Palta: add a new week of price data (2021-09-09, serial 44448)
# by inserting 7 new rows before the existing row 677, shifting the rest of the
# table down (old rows 677-684 become 684-691).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows above current row 677 (rows 677:683), pushing existing
# data (old rows 677-684) down to 684-691.
$ws.Range("A677:A683").EntireRow.Insert()

# Columns A, B, C, E, F, G, H, I, J are constant for every data row in this sheet.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100106
$producto = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"

$fecha = 44448

$rows = @(
    @{ Row = 677; Variedad = "Fuerte";            Calidad = "Primera";  Volumen = 50; PMin = 2400;  PMax = 2400;  PProm = 2400;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Quillota"; PrecioKilo = 2400; Factor = 1 },
    @{ Row = 678; Variedad = "Hass";               Calidad = "Especial"; Volumen = 48; PMin = 27000; PMax = 27000; PProm = 27000; Unidad = "`$/bandeja 10 kilos";              Origen = "Perú";                    PrecioKilo = 2700; Factor = 10 },
    @{ Row = 679; Variedad = "Hass";               Calidad = "Primera";  Volumen = 56; PMin = 24000; PMax = 24000; PProm = 24000; Unidad = "`$/bandeja 10 kilos";              Origen = "Perú";                    PrecioKilo = 2400; Factor = 10 },
    @{ Row = 680; Variedad = "Hass";               Calidad = "Segunda";  Volumen = 58; PMin = 21000; PMax = 21000; PProm = 21000; Unidad = "`$/bandeja 10 kilos";              Origen = "Perú";                    PrecioKilo = 2100; Factor = 10 },
    @{ Row = 681; Variedad = "Negra de La Cruz";   Calidad = "Primera";  Volumen = 35; PMin = 2000;  PMax = 2000;  PProm = 2000;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Quillota"; PrecioKilo = 2000; Factor = 1 },
    @{ Row = 682; Variedad = "Negra de La Cruz";   Calidad = "Segunda";  Volumen = 36; PMin = 1800;  PMax = 1800;  PProm = 1800;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Quillota"; PrecioKilo = 1800; Factor = 1 },
    @{ Row = 683; Variedad = "Negra de La Cruz";   Calidad = "Tercera";  Volumen = 30; PMin = 1500;  PMax = 1500;  PProm = 1500;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Quillota"; PrecioKilo = 1500; Factor = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKilo
    $ws.Cells.Item($row, 20).Value = $r.Factor

    # Match the existing date-cell formatting used throughout column D.
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
